$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(12,13,14,15,16,17,18,19,20,21,49,87,140,141,151,184,185,250)
foreach ($r in $rows) {
    $ws.Range("H$r").Value = "Y"
    $ws.Range("I$r").Value = "Y"
    $ws.Range("K$r").Value = "Y"
    $ws.Range("M$r").Value = "Y"
}

$ws.Range("N12").Select()
